$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date label used in column A (rows 2-7), stored as plain text.
# Setting NumberFormat to Text ("@") first stops Excel from auto-converting
# the literal "2025-12-01" into a date serial number when it's assigned.
$ws.Range("A2:A7").NumberFormat = "@"
$ws.Range("A2:A7").Value = "2025-12-01"
# Revert the cells back to the workbook's default (Normal) style so no
# extra formatting is introduced - only the text content itself changes.
$ws.Range("A2:A7").Style = "Normal"

# Refresh the computed probability figure in column N (rows 2-7).
$ws.Range("N2:N7").Value = 85.87246918135976
